$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "section" value (Table1[section], cell H3) ---
# This single cell drives the formulas in D2, C3, C9, C15, C21, C27, C33
# via the structured table reference Table1[section].
$ws.Range("H3").Value = "emptyHouseInterior"

# --- Update B column (navigation grid row index) for rows 3-38: +80 ---
for ($r = 3; $r -le 8; $r++)  { $ws.Range("B$r").Value = 128 }
for ($r = 9; $r -le 14; $r++) { $ws.Range("B$r").Value = 129 }
for ($r = 15; $r -le 20; $r++) { $ws.Range("B$r").Value = 130 }
for ($r = 21; $r -le 26; $r++) { $ws.Range("B$r").Value = 131 }
for ($r = 27; $r -le 32; $r++) { $ws.Range("B$r").Value = 132 }
for ($r = 33; $r -le 38; $r++) { $ws.Range("B$r").Value = 133 }

# --- Update C column text values (plain values, not the formula cells) ---
# Choices block (rows 4-6, row 3 is the formula-driven header)
$ws.Range("C4").Value = '"Take a look in the kitchen"'
$ws.Range("C5").Value = '"Check out the basement"'
$ws.Range("C6").Value = '"Head up to the bedroom"'

# Targets block (rows 10-12, row 9 is the formula-driven header)
$ws.Range("C10").Value = '"emptyHouseKitchen"'
$ws.Range("C11").Value = '"emptyHouseBasement"'
$ws.Range("C12").Value = '"emptyHouseBedroom"'

# InventoryGet block (rows 16-17 revert to NOTHING, row 15 is the formula-driven header)
$ws.Range("C16").Value = '"NOTHING"'
$ws.Range("C17").Value = '"NOTHING"'

# InventoryNeed block (row 23 gains a value, row 21 is the formula-driven header)
$ws.Range("C23").Value = '"First Aid Kit"'

# InventoryMissingTarget block (row 29 gains a value, row 27 is the formula-driven header)
$ws.Range("C29").Value = '"missingAidKit"'

# InventoryDestroy block (row 35 gains a value, row 33 is the formula-driven header)
$ws.Range("C35").Value = '"First Aid Kit"'

# --- Recalculate so dependent formulas (D column, C3/C9/C15/C21/C27/C33) refresh ---
$excel.Calculate()

# --- Column C width (target stored width is 46) ---
$ws.Columns("C").ColumnWidth = 45.1667

# --- Selection on the sheet ---
$ws.Range("L8").Select() | Out-Null

# --- Window vertical position (best effort; yWindow 1200 -> 1800) ---
$wb.Windows.Item(1).Top = 1800
